# The presentation currently uses the "Integral" theme (colour scheme) for
# its slide master / slides. The edit re-colours the deck's theme back to
# the default Office colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink),
# i.e. the equivalent of picking the "Office Theme" swatch from the Design
# gallery. Font scheme / format scheme are already identical between the two
# themes, so only the 12 theme colours need to change.

function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme, in ThemeColorScheme.Colors() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = HexToComRgb($officeThemeColors[$i - 1])
}
